$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.354576989141492
$ws.Cells.Item(2, 3).Value = 0.05622327177400166
$ws.Cells.Item(2, 4).Value = 0.3638058693529729
$ws.Cells.Item(2, 6).Value = 1.100482960602299
$ws.Cells.Item(2, 7).Value = 0.4918582983625797
$ws.Cells.Item(2, 8).Value = 0.6498530672995457
$ws.Cells.Item(2, 9).Value = 0.4978364578830146
$ws.Cells.Item(2, 10).Value = 0.3482663690695773
$ws.Cells.Item(2, 11).Value = 0.3746607540965101
$ws.Cells.Item(2, 15).Value = 2.24597801046886

$ws.Cells.Item(3, 2).Value = 0.3122309058202291
$ws.Cells.Item(3, 3).Value = 0.04942440441141116
$ws.Cells.Item(3, 4).Value = 0.3508189257738934
$ws.Cells.Item(3, 6).Value = 1.104065987383215
$ws.Cells.Item(3, 7).Value = 0.4966008789852268
$ws.Cells.Item(3, 8).Value = 0.6557557282262181
$ws.Cells.Item(3, 9).Value = 0.5047284444928302
$ws.Cells.Item(3, 10).Value = 0.3367206368869091
$ws.Cells.Item(3, 11).Value = 0.3280373710746289
$ws.Cells.Item(3, 15).Value = 2.268149801269601

$ws.Cells.Item(4, 2).Value = 0.2861698261037589
$ws.Cells.Item(4, 3).Value = 0.04522923779958887
$ws.Cells.Item(4, 4).Value = 0.342995951910666
$ws.Cells.Item(4, 6).Value = 1.106917768764816
$ws.Cells.Item(4, 7).Value = 0.499873451737102
$ws.Cells.Item(4, 8).Value = 0.659668871537427
$ws.Cells.Item(4, 9).Value = 0.5092540847112179
$ws.Cells.Item(4, 10).Value = 0.3298320962886976
$ws.Cells.Item(4, 11).Value = 0.2993212961998211
$ws.Cells.Item(4, 15).Value = 2.283124416136715

$ws.Cells.Item(5, 2).Value = 0.2755352827400088
$ws.Cells.Item(5, 3).Value = 0.0435145953079541
$ws.Cells.Item(5, 4).Value = 0.3398462065774339
$ws.Cells.Item(5, 6).Value = 1.108243808605337
$ws.Cells.Item(5, 7).Value = 0.5012976340262014
$ws.Cells.Item(5, 8).Value = 0.6613361837891603
$ws.Cells.Item(5, 9).Value = 0.5111721900710808
$ws.Cells.Item(5, 10).Value = 0.327075422774243
$ws.Cells.Item(5, 11).Value = 0.2875976534148492
$ws.Cells.Item(5, 15).Value = 2.289568797625904

$ws.Cells.Item(6, 2).Value = 0.2737685747786429
$ws.Cells.Item(6, 3).Value = 0.04322957679443107
$ws.Cells.Item(6, 4).Value = 0.3393255054626536
$ws.Cells.Item(6, 6).Value = 1.108473897160145
$ws.Cells.Item(6, 7).Value = 0.5015395865463432
$ws.Cells.Item(6, 8).Value = 0.6616174300261548
$ws.Cells.Item(6, 9).Value = 0.5114951503668976
$ws.Cells.Item(6, 10).Value = 0.3266207277699351
$ws.Cells.Item(6, 11).Value = 0.2856496700153173
$ws.Cells.Item(6, 15).Value = 2.290659540368949

$ws.Cells.Item(7, 2).Value = 0.2860264624111721
$ws.Cells.Item(7, 3).Value = 0.04520613394456063
$ws.Cells.Item(7, 4).Value = 0.3429533184972797
$ws.Cells.Item(7, 6).Value = 1.106934988489606
$ws.Cells.Item(7, 7).Value = 0.4998922921204567
$ws.Cells.Item(7, 8).Value = 0.6596910631949697
$ws.Cells.Item(7, 9).Value = 0.5092796539083455
$ws.Cells.Item(7, 10).Value = 0.3297947144717881
$ws.Cells.Item(7, 11).Value = 0.2991632734809855
$ws.Cells.Item(7, 15).Value = 2.283209942340875

$ws.Cells.Item(8, 2).Value = 0.3399890345640699
$ws.Cells.Item(8, 3).Value = 0.05388336929770787
$ws.Cells.Item(8, 4).Value = 0.3592967301938188
$ws.Cells.Item(8, 6).Value = 1.101583118434711
$ws.Cells.Item(8, 7).Value = 0.4934186380576335
$ws.Cells.Item(8, 8).Value = 0.6518283659384991
$ws.Cells.Item(8, 9).Value = 0.500151784823224
$ws.Cells.Item(8, 10).Value = 0.3442437791029676
$ws.Cells.Item(8, 11).Value = 0.3586039774162657
$ws.Cells.Item(8, 15).Value = 2.253340273830275

$ws.Cells.Item(9, 2).Value = 0.4453022527637245
$ws.Cells.Item(9, 3).Value = 0.07073160211504614
$ws.Cells.Item(9, 4).Value = 0.3925381009032947
$ws.Cells.Item(9, 6).Value = 1.096260078139693
$ws.Cells.Item(9, 7).Value = 0.4835895975822027
$ws.Cells.Item(9, 8).Value = 0.6387003923760375
$ws.Cells.Item(9, 9).Value = 0.4845859849575209
$ws.Cells.Item(9, 10).Value = 0.3741707445194749
$ws.Cells.Item(9, 11).Value = 0.4744300583580525
$ws.Cells.Item(9, 15).Value = 2.205572046147267

$ws.Cells.Item(10, 2).Value = 0.52233665861894
$ws.Cells.Item(10, 3).Value = 0.08300339220755859
$ws.Cells.Item(10, 4).Value = 0.4176811363624324
$ws.Cells.Item(10, 6).Value = 1.095504279907153
$ws.Cells.Item(10, 7).Value = 0.4781217560015918
$ws.Cells.Item(10, 8).Value = 0.6304498247039589
$ws.Cells.Item(10, 9).Value = 0.4745747218685779
$ws.Cells.Item(10, 10).Value = 0.3971327650132963
$ws.Cells.Item(10, 11).Value = 0.5590463496054383
$ws.Cells.Item(10, 15).Value = 2.17707447151983

$ws.Cells.Item(11, 2).Value = 0.5573020659880967
$ws.Cells.Item(11, 3).Value = 0.08856215963561453
$ws.Cells.Item(11, 4).Value = 0.4292747047656178
$ws.Cells.Item(11, 6).Value = 1.095846084125512
$ws.Cells.Item(11, 7).Value = 0.4760164721412181
$ws.Cells.Item(11, 8).Value = 0.6269988396135631
$ws.Cells.Item(11, 9).Value = 0.4703301962471578
$ws.Cells.Item(11, 10).Value = 0.4077914166588812
$ws.Cells.Item(11, 11).Value = 0.5974297861740752
$ws.Cells.Item(11, 15).Value = 2.165545047829866

$ws.Cells.Item(12, 2).Value = 0.5705306914873916
$ws.Cells.Item(12, 3).Value = 0.09066361141790935
$ws.Cells.Item(12, 4).Value = 0.4336871432487612
$ws.Cells.Item(12, 6).Value = 1.096074123298308
$ws.Cells.Item(12, 7).Value = 0.4752743076454991
$ws.Cells.Item(12, 8).Value = 0.6257354760532863
$ws.Cells.Item(12, 9).Value = 0.4687674783832456
$ws.Cells.Item(12, 10).Value = 0.4118582332731648
$ws.Cells.Item(12, 11).Value = 0.6119482263961231
$ws.Cells.Item(12, 15).Value = 2.161385585849345

$ws.Cells.Item(13, 2).Value = 0.5676822175381631
$ws.Cells.Item(13, 3).Value = 0.09021118520865912
$ws.Cells.Item(13, 4).Value = 0.4327358614002037
$ws.Cells.Item(13, 6).Value = 1.096020625275287
$ws.Cells.Item(13, 7).Value = 0.4754316952975941
$ws.Cells.Item(13, 8).Value = 0.6260056317586233
$ws.Cells.Item(13, 9).Value = 0.4691020533699373
$ws.Cells.Item(13, 10).Value = 0.4109810103449973
$ws.Cells.Item(13, 11).Value = 0.6088221696903418
$ws.Cells.Item(13, 15).Value = 2.162272213410091

$ws.Cells.Item(14, 2).Value = 0.5583906384011357
$ws.Cells.Item(14, 3).Value = 0.08873511882728735
$ws.Cells.Item(14, 4).Value = 0.4296372748985107
$ws.Cells.Item(14, 6).Value = 1.095862868772471
$ws.Cells.Item(14, 7).Value = 0.4759543096345524
$ws.Cells.Item(14, 8).Value = 0.626894031206767
$ws.Cells.Item(14, 9).Value = 0.4702007364496801
$ws.Cells.Item(14, 10).Value = 0.4081253826776816
$ws.Cells.Item(14, 11).Value = 0.5986245633244778
$ws.Cells.Item(14, 15).Value = 2.165198706872246

$ws.Cells.Item(15, 2).Value = 0.5526976891626703
$ws.Cells.Item(15, 3).Value = 0.0878305219135882
$ws.Cells.Item(15, 4).Value = 0.4277421868283113
$ws.Cells.Item(15, 6).Value = 1.095779080076355
$ws.Cells.Item(15, 7).Value = 0.4762816001055157
$ws.Cells.Item(15, 8).Value = 0.6274438595850427
$ws.Cells.Item(15, 9).Value = 0.4708795206591248
$ws.Cells.Item(15, 10).Value = 0.4063802145781921
$ws.Cells.Item(15, 11).Value = 0.5923760560027915
$ws.Cells.Item(15, 15).Value = 2.167018165422718

$ws.Cells.Item(16, 2).Value = 0.5200499496259852
$ws.Cells.Item(16, 3).Value = 0.08263962585486695
$ws.Cells.Item(16, 4).Value = 0.4169265881650404
$ws.Cells.Item(16, 6).Value = 1.095495739855878
$ws.Cells.Item(16, 7).Value = 0.4782670396825424
$ws.Cells.Item(16, 8).Value = 0.6306814354917023
$ws.Cells.Item(16, 9).Value = 0.4748583480597546
$ws.Cells.Item(16, 10).Value = 0.3964404830072823
$ws.Cells.Item(16, 11).Value = 0.5565356352804542
$ws.Cells.Item(16, 15).Value = 2.177856827684579

$ws.Cells.Item(17, 2).Value = 0.5000010824323908
$ws.Cells.Item(17, 3).Value = 0.07944901535294946
$ws.Cells.Item(17, 4).Value = 0.4103313436432927
$ws.Cells.Item(17, 6).Value = 1.095497533997872
$ws.Cells.Item(17, 7).Value = 0.479582984435666
$ws.Cells.Item(17, 8).Value = 0.632744985591799
$ws.Cells.Item(17, 9).Value = 0.4773785843109444
$ws.Cells.Item(17, 10).Value = 0.3903973429306831
$ws.Cells.Item(17, 11).Value = 0.5345201769616779
$ws.Cells.Item(17, 15).Value = 2.18487351347126

$ws.Cells.Item(18, 2).Value = 0.4884622173137529
$ws.Cells.Item(18, 3).Value = 0.07761163402373938
$ws.Cells.Item(18, 4).Value = 0.4065526176958656
$ws.Cells.Item(18, 6).Value = 1.09556309904287
$ws.Cells.Item(18, 7).Value = 0.4803758356517989
$ws.Cells.Item(18, 8).Value = 0.6339603339546684
$ws.Cells.Item(18, 9).Value = 0.4788572992298263
$ws.Cells.Item(18, 10).Value = 0.3869415496426711
$ws.Cells.Item(18, 11).Value = 0.5218472655326991
$ws.Cells.Item(18, 15).Value = 2.189044311509846

$ws.Cells.Item(19, 2).Value = 0.4845541290070798
$ws.Cells.Item(19, 3).Value = 0.07698914953422786
$ws.Cells.Item(19, 4).Value = 0.4052757344114752
$ws.Cells.Item(19, 6).Value = 1.095596381420272
$ws.Cells.Item(19, 7).Value = 0.4806504530053033
$ws.Cells.Item(19, 8).Value = 0.634376716529907
$ws.Cells.Item(19, 9).Value = 0.4793629698281947
$ws.Cells.Item(19, 10).Value = 0.3857749237761965
$ws.Cells.Item(19, 11).Value = 0.5175547097831839
$ws.Cells.Item(19, 15).Value = 2.190479648850655

$ws.Cells.Item(20, 2).Value = 0.5021360786163882
$ws.Cells.Item(20, 3).Value = 0.07978889260138544
$ws.Cells.Item(20, 4).Value = 0.4110319004278153
$ws.Cells.Item(20, 6).Value = 1.09549066396913
$ws.Cells.Item(20, 7).Value = 0.4794391778288585
$ws.Cells.Item(20, 8).Value = 0.6325223728912803
$ws.Cells.Item(20, 9).Value = 0.4771072843340995
$ws.Cells.Item(20, 10).Value = 0.3910385694715757
$ws.Cells.Item(20, 11).Value = 0.5368648222653576
$ws.Cells.Item(20, 15).Value = 2.184112603042351

$ws.Cells.Item(21, 2).Value = 0.5611201329543292
$ws.Cells.Item(21, 3).Value = 0.08916877213948737
$ws.Cells.Item(21, 4).Value = 0.4305468038758136
$ws.Cells.Item(21, 6).Value = 1.095906529456215
$ws.Cells.Item(21, 7).Value = 0.4757993098459679
$ws.Cells.Item(21, 8).Value = 0.6266319076395703
$ws.Cells.Item(21, 9).Value = 0.4698768158702826
$ws.Cells.Item(21, 10).Value = 0.4089633190565962
$ws.Cells.Item(21, 11).Value = 0.6016203031500993
$ws.Cells.Item(21, 15).Value = 2.16433351911725

$ws.Cells.Item(22, 2).Value = 0.5995992136623158
$ws.Cells.Item(22, 3).Value = 0.09527841832198192
$ws.Cells.Item(22, 4).Value = 0.4434302560893002
$ws.Cells.Item(22, 6).Value = 1.096753074516968
$ws.Cells.Item(22, 7).Value = 0.4737414374487514
$ws.Cells.Item(22, 8).Value = 0.623035400100342
$ws.Cells.Item(22, 9).Value = 0.4654112322640032
$ws.Cells.Item(22, 10).Value = 0.4208566215411196
$ws.Cells.Item(22, 11).Value = 0.6438449784263867
$ws.Cells.Item(22, 15).Value = 2.152610378831426

$ws.Cells.Item(23, 2).Value = 0.5790689131097224
$ws.Cells.Item(23, 3).Value = 0.09201951400999064
$ws.Cells.Item(23, 4).Value = 0.4365423502285637
$ws.Cells.Item(23, 6).Value = 1.09624866327772
$ws.Cells.Item(23, 7).Value = 0.4748103523823346
$ws.Cells.Item(23, 8).Value = 0.6249317552324243
$ws.Cells.Item(23, 9).Value = 0.4677707918667977
$ws.Cells.Item(23, 10).Value = 0.414492624547961
$ws.Cells.Item(23, 11).Value = 0.6213180157621991
$ws.Cells.Item(23, 15).Value = 2.158757029641265

$ws.Cells.Item(24, 2).Value = 0.5011708857298629
$ws.Cells.Item(24, 3).Value = 0.07963524361019836
$ws.Cells.Item(24, 4).Value = 0.4107151382963252
$ws.Cells.Item(24, 6).Value = 1.095493568892763
$ws.Cells.Item(24, 7).Value = 0.4795040797422061
$ws.Cells.Item(24, 8).Value = 0.6326229258021172
$ws.Cells.Item(24, 9).Value = 0.4772298462411193
$ws.Cells.Item(24, 10).Value = 0.3907486133628879
$ws.Cells.Item(24, 11).Value = 0.5358048576820522
$ws.Cells.Item(24, 15).Value = 2.184456184504285

$ws.Cells.Item(25, 2).Value = 0.416869651338402
$ws.Cells.Item(25, 3).Value = 0.06619212042424749
$ws.Cells.Item(25, 4).Value = 0.3834184442833646
$ws.Cells.Item(25, 6).Value = 1.097146189578567
$ws.Cells.Item(25, 7).Value = 0.4859411376194984
$ws.Cells.Item(25, 8).Value = 0.642006818713007
$ws.Cells.Item(25, 9).Value = 0.4885468655629701
$ws.Cells.Item(25, 10).Value = 0.3659038219808366
$ws.Cells.Item(25, 11).Value = 0.4431782373231101
$ws.Cells.Item(25, 15).Value = 2.217336706370588

Write-Host "done"